$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Stikine" regional district block: its header/summary row,
# its single "Unincorporated Areas" subarea row, and the blank separator
# row beneath it. Deleting with a shift pulls every subsequent row (the
# "Sunshine Coast" section onward) up by three rows, which also drops the
# now-unused "57000" / "Stikine" / "R" / "57999" shared strings.
$ws.Rows("226:228").Delete()

# Restore the cursor/selection to where it was left in the saved workbook.
$ws.Range("B216").Select()
